$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(44462, 800, 9000, 10000, 9500, 380)
    3 = @(44349, 600, 10000, 12000, 11000, 440)
    4 = @(44406, 800, 10000, 11000, 10500, 420)
    5 = @(44419, 1100, 11000, 12000, 11500, 460)
    6 = @(44475, 1200, 5000, 6000, 5500, 220)
    7 = @(44363, 900, 11000, 12000, 11500, 460)
    8 = @(44391, 500, 9000, 10000, 9500, 380)
    9 = @(44364, 700, 11000, 12000, 11500, 460)
    10 = @(44434, 600, 10000, 11000, 10500, 420)
    11 = @(44441, 1100, 11000, 12000, 11500, 460)
    12 = @(44413, 1200, 10000, 11000, 10500, 420)
    13 = @(44448, 800, 10000, 12000, 11000, 440)
    14 = @(44426, 500, 11000, 12000, 11500, 460)
    15 = @(44489, 1200, 5000, 6000, 5500, 220)
    16 = @(44461, 1100, 9000, 10000, 9500, 380)
    17 = @(44455, 600, 9000, 10000, 9500, 380)
    18 = @(44447, 1000, 10000, 12000, 11000, 440)
    19 = @(44483, 1200, 4000, 5000, 4500, 180)
    20 = @(44435, 600, 10000, 11000, 10500, 420)
    21 = @(44427, 360, 10000, 11000, 10500, 420)
    22 = @(44468, 700, 5000, 6000, 5500, 220)
    23 = @(44336, 1200, 12000, 13000, 12500, 500)
    24 = @(44420, 1000, 10000, 11000, 10500, 420)
    25 = @(44412, 1000, 10000, 11000, 10500, 420)
    26 = @(44377, 800, 9000, 10000, 9500, 380)
    27 = @(44335, 1000, 12000, 13000, 12500, 500)
    28 = @(44385, 600, 8000, 9000, 8500, 340)
    29 = @(44308, 400, 11000, 12000, 11500, 460)
    30 = @(44398, 400, 9000, 10000, 9500, 380)
    31 = @(44371, 500, 10000, 12000, 11000, 440)
    32 = @(44454, 800, 9000, 10000, 9500, 380)
    33 = @(44490, 400, 5000, 6000, 5500, 220)
    34 = @(44399, 500, 9000, 10000, 9500, 380)
    35 = @(44329, 1000, 12000, 13000, 12500, 500)
    36 = @(44476, 1100, 5000, 6000, 5500, 220)
    37 = @(44482, 1600, 4000, 5000, 4500, 180)
    38 = @(44356, 1000, 11000, 12000, 11500, 460)
    39 = @(44469, 600, 5000, 6000, 5500, 220)
    40 = @(44384, 700, 8000, 9000, 8500, 340)
    41 = @(44343, 500, 9000, 10000, 9500, 380)
    42 = @(44392, 600, 9000, 10000, 9500, 380)
    43 = @(44328, 900, 11000, 12000, 11500, 460)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
    $ws.Cells.Item($row, 11).Value = $vals[2]
    $ws.Cells.Item($row, 12).Value = $vals[3]
    $ws.Cells.Item($row, 13).Value = $vals[4]
    $ws.Cells.Item($row, 16).Value = $vals[5]
}

Write-Host "Done updating rows."